$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.014.29"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "2.052.07"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.64"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.77"
$ws.Range("E7").Value = "  +6.49%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("D12").Value = "2.355.09"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.65"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.751"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "2.068.21"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "37.925.12"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.71"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.77"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.34"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.01"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("E34").Value = "  +10.30%  "
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.33"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  +9.92%  "
$ws.Range("E38").Value = "  +4.97%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "1.486.98"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.15"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.54"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.13"
$ws.Range("E47").Value = "  +12.21%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.07"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "2.244.91"
$ws.Range("E51").Value = "  +1.71%  "
